$d = $word.ActiveDocument

# Locate the run containing " Classic & Quantum Mechanics" and replace just
# that run's text with " Material" (ReplaceOne keeps the match a single run
# and repositions the range onto the replacement).
$rng = $d.Content
$rng.Find.Execute(" Classic & Quantum Mechanics", $true, $false, $false, $false, $false,
                   $true, 1, $false, " Material", 1)

# Move just past the replacement, then replace the very next single-space
# run (the one that sits before "Engineer") with " Scientists & ", leaving
# the "Engineer" run untouched.
$rng.Collapse(0)
$rng.Find.Execute(" ", $false, $false, $false, $false, $false,
                   $true, 1, $false, " Scientists & ", 1)
